# Fruta / hortaliza, semanal
# Insert 5 new "Edranol" price rows (week of 2021-09-29, serial 44468) right
# before the existing row 182 ("Comercializadora del Agro de Limarí" / Palta
# block), pushing the rest of the table down by 5 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 5 blank rows at 182..186 (rows below shift down).
$ws.Rows("182:186").Insert()

$newRows = @(
    @(44468, "Edranol", "Especial",     200, 2300, 2400, 2350, "$/kilo (en caja de 17 kilos)", 2350),
    @(44468, "Edranol", "Primera",      240, 2100, 2200, 2150, "$/kilo (en caja de 17 kilos)", 2150),
    @(44468, "Edranol", "Segunda",      240, 1800, 1900, 1850, "$/kilo (en caja de 17 kilos)", 1850),
    @(44468, "Edranol", "1a nueva(o)",  500, 2300, 2400, 2350, "$/kilo (en caja de 17 kilos)", 2350),
    @(44468, "Edranol", "2a nueva(o)",  400, 2100, 2200, 2150, "$/kilo (en caja de 17 kilos)", 2150)
)

$r = 182
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = 2
    $ws.Cells.Item($r, 2).Value  = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value  = "Coquimbo"
    $ws.Cells.Item($r, 4).Value  = $row[0]
    $ws.Cells.Item($r, 5).Value  = 4
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100106
    $ws.Cells.Item($r, 8).Value  = "Oleaginosos"
    $ws.Cells.Item($r, 9).Value  = 100106002
    $ws.Cells.Item($r, 10).Value = "Palta"
    $ws.Cells.Item($r, 11).Value = $row[1]
    $ws.Cells.Item($r, 12).Value = $row[2]
    $ws.Cells.Item($r, 13).Value = $row[3]
    $ws.Cells.Item($r, 14).Value = $row[4]
    $ws.Cells.Item($r, 15).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 17).Value = $row[7]
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 19).Value = $row[8]
    $ws.Cells.Item($r, 20).Value = 1
    $r = $r + 1
}
